# Apply edits to the "Deviations" worksheet (sheet3.xml):
#  - F2: 1 -> 1.2
#  - A3: 25 -> 120
#  - Delete rows 4 through 10 (data that used to follow row 3)
#  - Update selection to F5 (sqref A4:F5)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deviations")

# Update individual cell values first
$ws.Range("F2").Value = 1.2
$ws.Range("A3").Value = 120
$ws.Range("F3").Value = 1

# Delete rows 4 through 10 entirely (shifts rows up, shrinking used range)
$ws.Range("A4:F10").EntireRow.Delete()

# Update the active selection / active cell as recorded in the saved file
$ws.Activate()
$ws.Range("A4:F5").Select()
$excel.ActiveCell = $ws.Range("F5")
